$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of data (row 22) - values reuse some already-existing strings
# and introduce three new ones (version, folder, changelog).
$ws.Range("A22").Value = "0.6.3"
$ws.Range("B22").Value = "AUTOMATA CELULAR - copia (32)"
$ws.Range("C22").Value = $ws.Range("C20").Value2
$ws.Range("D22").Value = "-Reworked Greed calc`n-GUI: Flexibility has now the 0 value by default`n-Fixed Flexibility not working properly."
$ws.Range("E22").Value = $ws.Range("E20").Value2
$ws.Range("F22").Value = $ws.Range("F20").Value2

# Row height for the newly filled-in row matches rows 20/21 wrapped style.
$ws.Rows.Item(22).RowHeight = 57.6

# Move the active selection down one row (to D23) to reflect where the
# user ended up after entering the new data, and drop the prior
# horizontal scroll position (topLeftCell was D1, now back to default).
$ws.Range("D23").Select()
